# The deck's Design theme colours are switched from the custom
# "Integral" (Red Violet) scheme to the default "Office Theme" colour
# scheme: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink on the theme part
# used by the slide master (ppt/theme/theme1.xml) all change to the
# stock Office RGB values.
#
# In the PowerPoint object model theme colours are read/written
# through a slide's ThemeColorScheme (12 entries, in ppColorSchemeIndex
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) - this maps
# straight onto <a:clrScheme> in the theme XML shared by every slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$scheme = $s.ThemeColorScheme

function ToComRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" colour scheme values.
$scheme.Colors(1).RGB  = ToComRGB 0x00 0x00 0x00   # dk1       000000
$scheme.Colors(2).RGB  = ToComRGB 0xFF 0xFF 0xFF   # lt1       FFFFFF
$scheme.Colors(3).RGB  = ToComRGB 0x44 0x54 0x6A   # dk2       44546A
$scheme.Colors(4).RGB  = ToComRGB 0xE7 0xE6 0xE6   # lt2       E7E6E6
$scheme.Colors(5).RGB  = ToComRGB 0x5B 0x9B 0xD5   # accent1   5B9BD5
$scheme.Colors(6).RGB  = ToComRGB 0xED 0x7D 0x31   # accent2   ED7D31
$scheme.Colors(7).RGB  = ToComRGB 0xA5 0xA5 0xA5   # accent3   A5A5A5
$scheme.Colors(8).RGB  = ToComRGB 0xFF 0xC0 0x00   # accent4   FFC000
$scheme.Colors(9).RGB  = ToComRGB 0x44 0x72 0xC4   # accent5   4472C4
$scheme.Colors(10).RGB = ToComRGB 0x70 0xAD 0x47   # accent6   70AD47
$scheme.Colors(11).RGB = ToComRGB 0x05 0x63 0xC1   # hlink     0563C1
$scheme.Colors(12).RGB = ToComRGB 0x95 0x4F 0x72   # folHlink  954F72
